$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A118").Value = "80266890"
$ws.Range("B118").Value = "10050-ARI-I"
$ws.Range("C118").Value = 2
$ws.Range("A119").Value = "80266890"
$ws.Range("B119").Value = "10496-ARI-I"
$ws.Range("C119").Value = 2
$ws.Range("A120").Value = "80266890"
$ws.Range("B120").Value = "10070-ARI-I"
$ws.Range("C120").Value = 2
$ws.Range("A121").Value = "80266891"
$ws.Range("B121").Value = "60164-YAG-I"
$ws.Range("C121").Value = 5000
$ws.Range("A122").Value = "80266892"
$ws.Range("B122").Value = "40274-TDK-I"
$ws.Range("C122").Value = 200
$ws.Range("A123").Value = "80266892"
$ws.Range("B123").Value = "21487-MET-I"
$ws.Range("C123").Value = 1200
$ws.Range("A124").Value = "80266892"
$ws.Range("B124").Value = "40312-TDK-N"
$ws.Range("C124").Value = 500
$ws.Range("A125").Value = "80266893"
$ws.Range("B125").Value = "26489-YAG-I"
$ws.Range("C125").Value = 3000
$ws.Range("A126").Value = "80266893"
$ws.Range("B126").Value = "16512-YAG-I"
$ws.Range("C126").Value = 8000
$ws.Range("A127").Value = "80266893"
$ws.Range("B127").Value = "30056-EVL-I"
$ws.Range("C127").Value = 4000
$ws.Range("A128").Value = "80266893"
$ws.Range("B128").Value = "17498-VIS-I"
$ws.Range("C128").Value = 3200
$ws.Range("A129").Value = "80266893"
$ws.Range("B129").Value = "30287-ONS-L"
$ws.Range("C129").Value = 2500
$ws.Range("A130").Value = "80266893"
$ws.Range("B130").Value = "10136-ROY-I"
$ws.Range("C130").Value = 5000
$ws.Range("A131").Value = "80266893"
$ws.Range("B131").Value = "10587-ROY-I"
$ws.Range("C131").Value = 5000
$ws.Range("A132").Value = "80266893"
$ws.Range("B132").Value = "10917-ROY-I"
$ws.Range("C132").Value = 5000
$ws.Range("A133").Value = "80266893"
$ws.Range("B133").Value = "11587-ROY-I"
$ws.Range("C133").Value = 5000
$ws.Range("A134").Value = "80266893"
$ws.Range("B134").Value = "11875-ROY-I"
$ws.Range("C134").Value = 5000
$ws.Range("A135").Value = "80266893"
$ws.Range("B135").Value = "10001-THK-I"
$ws.Range("C135").Value = 40
$ws.Range("A136").Value = "80266893"
$ws.Range("B136").Value = "19880-WRN-I"
$ws.Range("C136").Value = 3000
$ws.Range("A137").Value = "80266893"
$ws.Range("B137").Value = "10000-THK-I"
$ws.Range("C137").Value = 40
$ws.Range("A138").Value = "80266893"
$ws.Range("B138").Value = "20606-EPC-I"
$ws.Range("C138").Value = 500
$ws.Range("A139").Value = "80266893"
$ws.Range("B139").Value = "60255-OUT-L"
$ws.Range("C139").Value = 350
$ws.Range("A140").Value = "80266893"
$ws.Range("B140").Value = "60260-OUT-L"
$ws.Range("C140").Value = 2200
$ws.Range("A141").Value = "80266893"
$ws.Range("B141").Value = "40312-TDK-N"
$ws.Range("C141").Value = 500
$ws.Range("A142").Value = "80266896"
$ws.Range("B142").Value = "10259-ARI-I"
$ws.Range("C142").Value = 18
$ws.Range("A143").Value = "80266898"
$ws.Range("B143").Value = "10259-ARI-I"
$ws.Range("C143").Value = 22
$ws.Range("A144").Value = "80266898"
$ws.Range("B144").Value = "10160-ARI-I"
$ws.Range("C144").Value = 4
$ws.Range("A145").Value = "80266899"
$ws.Range("B145").Value = "10259-ARI-I"
$ws.Range("C145").Value = 28
$ws.Range("A146").Value = "80266906"
$ws.Range("B146").Value = "10378-ARI-I"
$ws.Range("C146").Value = 12
$ws.Range("A147").Value = "80266906"
$ws.Range("B147").Value = "10258-ARI-I"
$ws.Range("C147").Value = 8
$ws.Range("A148").Value = "80266906"
$ws.Range("B148").Value = "10256-ARI-I"
$ws.Range("C148").Value = 6
$ws.Range("A149").Value = "80266906"
$ws.Range("B149").Value = "10493-ARI-I"
$ws.Range("C149").Value = 12
$ws.Range("A150").Value = "80266906"
$ws.Range("B150").Value = "10482-ARI-I"
$ws.Range("C150").Value = 2
$ws.Range("A151").Value = "80266907"
$ws.Range("B151").Value = "15395-DLO-I"
$ws.Range("C151").Value = 4
$ws.Range("A152").Value = "80266907"
$ws.Range("B152").Value = "11558-DLO-I"
$ws.Range("C152").Value = 1
$ws.Range("A153").Value = "80266907"
$ws.Range("B153").Value = "14755-DLO-I"
$ws.Range("C153").Value = 1
$ws.Range("A154").Value = "80266907"
$ws.Range("B154").Value = "10967-DLO-L"
$ws.Range("C154").Value = 2
$ws.Range("A155").Value = "80266907"
$ws.Range("B155").Value = "11465-DLO-I"
$ws.Range("C155").Value = 1
$ws.Range("A156").Value = "80266907"
$ws.Range("B156").Value = "15151-DLO-I"
$ws.Range("C156").Value = 1
$ws.Range("A157").Value = "80266907"
$ws.Range("B157").Value = "15141-DLO-I"
$ws.Range("C157").Value = 10
$ws.Range("A158").Value = "80266908"
$ws.Range("B158").Value = "20935-CTY-I"
$ws.Range("C158").Value = 1
$ws.Range("A159").Value = "80266909"
$ws.Range("B159").Value = "10494-ARI-I"
$ws.Range("C159").Value = 1
$ws.Range("A160").Value = "80266912"
$ws.Range("B160").Value = "10359-ARI-I"
$ws.Range("C160").Value = 2
$ws.Range("A161").Value = "80266912"
$ws.Range("B161").Value = "10645-ARI-I"
$ws.Range("C161").Value = 2
$ws.Range("A162").Value = "80266912"
$ws.Range("B162").Value = "10637-ARI-I"
$ws.Range("C162").Value = 2
$ws.Range("A163").Value = "80266912"
$ws.Range("B163").Value = "10488-ARI-I"
$ws.Range("C163").Value = 2
$ws.Range("A164").Value = "80266912"
$ws.Range("B164").Value = "10396-ARI-I"
$ws.Range("C164").Value = 2
$ws.Range("A165").Value = "80266913"
$ws.Range("B165").Value = "10652-ARI-I"
$ws.Range("C165").Value = 2
$ws.Range("A166").Value = "80266913"
$ws.Range("B166").Value = "10359-ARI-I"
$ws.Range("C166").Value = 4
$ws.Range("A167").Value = "80266913"
$ws.Range("B167").Value = "10035-ARI-I"
$ws.Range("C167").Value = 1
$ws.Range("A168").Value = "80266913"
$ws.Range("B168").Value = "10055-ARI-I"
$ws.Range("C168").Value = 2
$ws.Range("A169").Value = "80266913"
$ws.Range("B169").Value = "10547-ARI-I"
$ws.Range("C169").Value = 2
$ws.Range("A170").Value = "80266913"
$ws.Range("B170").Value = "10540-ARI-I"
$ws.Range("C170").Value = 2
$ws.Range("A171").Value = "80266913"
$ws.Range("B171").Value = "10645-ARI-I"
$ws.Range("C171").Value = 4
$ws.Range("A172").Value = "80266913"
$ws.Range("B172").Value = "10637-ARI-I"
$ws.Range("C172").Value = 4
$ws.Range("A173").Value = "80266913"
$ws.Range("B173").Value = "10488-ARI-I"
$ws.Range("C173").Value = 4
$ws.Range("A174").Value = "80266913"
$ws.Range("B174").Value = "10396-ARI-I"
$ws.Range("C174").Value = 4
$ws.Range("A175").Value = "80266917"
$ws.Range("B175").Value = "30259-OSR-I"
$ws.Range("C175").Value = 7500

$ws.Range("A2:C175").Select()
